$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '50.396.71'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -16.94%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.198.94'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -24.32%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '409.90'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -22.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.31'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -21.68%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.447'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -19.37%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.196.56'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -24.60%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '4.97'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -17.64%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0835'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -22.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.118'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -6.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.280'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -22.18%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.581.83'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -24.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '50.393.24'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -16.91%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '17.73'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -22.21%  '
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.225.49'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -23.71%  '
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000110'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -22.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.69'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -26.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '280.48'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -22.31%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.28'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -29.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.96'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -25.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '51.63'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -20.42%  '
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.341'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -24.94%  '
$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.134'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -25.96%  '
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'USDe'
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.47'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -17.75%  '
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0620'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -26.90%  '
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '138.41'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -7.00%  '
$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '16.10'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -18.76%  '
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.23'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -27.17%  '
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.42'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -20.80%  '
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.23'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -25.85%  '
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.753'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -25.04%  '
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.931'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -22.99%  '
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '10.19'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '30.72'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -18.96%  '
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -21.16%  '
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.13'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -24.10%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0464'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -20.87%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.783.47'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -22.33%  '
$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.493'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -24.13%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0782'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -15.29%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0191'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -20.09%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'ZEEBU'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.56'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -5.77%  '
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '14.56'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -28.99%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.60'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -28.09%  '
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.855'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -15.75%  '
